$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Support to skip special characters:
# update the SkipScenario (L2) and MultiRun (M2) sample values on the
# PetPost row, and move the active selection from L3 to M3.
$ws.Range("L2").Value = 'EXACT ("[petName]","Test2")'
$ws.Range("M2").Value = "petId|petName;1000|Test1;2000|Test2"

[void]$ws.Range("M3").Select()
